$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spelling corrections in the "Description" column (column B) of the metadata sheet.
# B8  = FIRE_ID description: "identifyier" -> "identifier"
$ws.Range("B8").Value = "fire incident identifier "

# B13 = heat_load description: "calcualted" -> "calculated"
$ws.Range("B13").Value = "heat load calculated at 30m resolution following McCune & Keon 2002"

# B15 = LYGrowth description: "calcualted" -> "calculated"
$ws.Range("B15").Value = "last year's growth, calculated as this year's height minus LYHeight (cm)"

# Update the view: scroll down and move the active cell selection.
$ws.Range("I36").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
